$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) for rows 2-3 ---
$sched = $wb.Worksheets.Item("Schedule")
$sched.Range("E2").Value = 1161.84723525
$sched.Range("F2").Value = 25.61391612103175
$sched.Range("E3").Value = 348.247965
$sched.Range("F3").Value = 23.03227281746032

# --- Sheet "Detailed": insert a new row at position 3 (shifts old rows 3-48 down to 4-49) ---
$ws = $wb.Worksheets.Item("Detailed")
$ws.Rows.Item(3).Insert()

# Now rewrite rows 2-49 with the final values (run 27 results)
$ws.Range("A2").Value = 46038
$ws.Range("B2").Value = 81.05549
$ws.Range("C2").Value = "historical"
$ws.Range("D2").Value = 46038
$ws.Range("E2").Value = "OFF"
$ws.Range("A3").Value = 46038.02083333334
$ws.Range("B3").Value = 78
$ws.Range("C3").Value = "historical"
$ws.Range("D3").Value = 46038
$ws.Range("E3").Value = "OFF"
$ws.Range("A4").Value = 46038.04166666666
$ws.Range("B4").Value = 78
$ws.Range("C4").Value = "historical"
$ws.Range("D4").Value = 46038
$ws.Range("E4").Value = "OFF"
$ws.Range("A5").Value = 46038.0625
$ws.Range("B5").Value = 84.7901
$ws.Range("C5").Value = "historical"
$ws.Range("D5").Value = 46038
$ws.Range("E5").Value = "OFF"
$ws.Range("A6").Value = 46038.08333333334
$ws.Range("B6").Value = 84.7901
$ws.Range("C6").Value = "historical"
$ws.Range("D6").Value = 46038
$ws.Range("E6").Value = "OFF"
$ws.Range("A7").Value = 46038.10416666666
$ws.Range("B7").Value = 83.60652
$ws.Range("C7").Value = "forecast"
$ws.Range("D7").Value = 46038
$ws.Range("E7").Value = "OFF"
$ws.Range("A8").Value = 46038.125
$ws.Range("B8").Value = 84.7901
$ws.Range("C8").Value = "forecast"
$ws.Range("D8").Value = 46038
$ws.Range("E8").Value = "OFF"
$ws.Range("A9").Value = 46038.14583333334
$ws.Range("B9").Value = 84.7901
$ws.Range("C9").Value = "forecast"
$ws.Range("D9").Value = 46038
$ws.Range("E9").Value = "OFF"
$ws.Range("A10").Value = 46038.16666666666
$ws.Range("B10").Value = 84.7901
$ws.Range("C10").Value = "forecast"
$ws.Range("D10").Value = 46038
$ws.Range("E10").Value = "OFF"
$ws.Range("A11").Value = 46038.1875
$ws.Range("B11").Value = 77.94
$ws.Range("C11").Value = "forecast"
$ws.Range("D11").Value = 46038
$ws.Range("E11").Value = "ON"
$ws.Range("A12").Value = 46038.20833333334
$ws.Range("B12").Value = 79.95
$ws.Range("C12").Value = "forecast"
$ws.Range("D12").Value = 46038
$ws.Range("E12").Value = "ON"
$ws.Range("A13").Value = 46038.22916666666
$ws.Range("B13").Value = 80.02
$ws.Range("C13").Value = "forecast"
$ws.Range("D13").Value = 46038
$ws.Range("E13").Value = "ON"
$ws.Range("A14").Value = 46038.25
$ws.Range("B14").Value = 79.95
$ws.Range("C14").Value = "forecast"
$ws.Range("D14").Value = 46038
$ws.Range("E14").Value = "ON"
$ws.Range("A15").Value = 46038.27083333334
$ws.Range("B15").Value = 57.06003
$ws.Range("C15").Value = "forecast"
$ws.Range("D15").Value = 46038
$ws.Range("E15").Value = "ON"
$ws.Range("A16").Value = 46038.29166666666
$ws.Range("B16").Value = 50.63034
$ws.Range("C16").Value = "forecast"
$ws.Range("D16").Value = 46038
$ws.Range("E16").Value = "ON"
$ws.Range("A17").Value = 46038.3125
$ws.Range("B17").Value = 56.98
$ws.Range("C17").Value = "forecast"
$ws.Range("D17").Value = 46038
$ws.Range("E17").Value = "ON"
$ws.Range("A18").Value = 46038.33333333334
$ws.Range("B18").Value = 56.98
$ws.Range("C18").Value = "forecast"
$ws.Range("D18").Value = 46038
$ws.Range("E18").Value = "ON"
$ws.Range("A19").Value = 46038.35416666666
$ws.Range("B19").Value = 56.89739
$ws.Range("C19").Value = "forecast"
$ws.Range("D19").Value = 46038
$ws.Range("E19").Value = "ON"
$ws.Range("A20").Value = 46038.375
$ws.Range("B20").Value = 46.64753
$ws.Range("C20").Value = "forecast"
$ws.Range("D20").Value = 46038
$ws.Range("E20").Value = "ON"
$ws.Range("A21").Value = 46038.39583333334
$ws.Range("B21").Value = 36.06
$ws.Range("C21").Value = "forecast"
$ws.Range("D21").Value = 46038
$ws.Range("E21").Value = "ON"
$ws.Range("A22").Value = 46038.41666666666
$ws.Range("B22").Value = 46.76441
$ws.Range("C22").Value = "forecast"
$ws.Range("D22").Value = 46038
$ws.Range("E22").Value = "ON"
$ws.Range("A23").Value = 46038.4375
$ws.Range("B23").Value = 41.46457
$ws.Range("C23").Value = "forecast"
$ws.Range("D23").Value = 46038
$ws.Range("E23").Value = "ON"
$ws.Range("A24").Value = 46038.45833333334
$ws.Range("B24").Value = 36.06028
$ws.Range("C24").Value = "forecast"
$ws.Range("D24").Value = 46038
$ws.Range("E24").Value = "ON"
$ws.Range("A25").Value = 46038.47916666666
$ws.Range("B25").Value = 36.06
$ws.Range("C25").Value = "forecast"
$ws.Range("D25").Value = 46038
$ws.Range("E25").Value = "ON"
$ws.Range("A26").Value = 46038.5
$ws.Range("B26").Value = 36.06
$ws.Range("C26").Value = "forecast"
$ws.Range("D26").Value = 46038
$ws.Range("E26").Value = "ON"
$ws.Range("A27").Value = 46038.52083333334
$ws.Range("B27").Value = 41.10377
$ws.Range("C27").Value = "forecast"
$ws.Range("D27").Value = 46038
$ws.Range("E27").Value = "ON"
$ws.Range("A28").Value = 46038.54166666666
$ws.Range("B28").Value = 36.0601
$ws.Range("C28").Value = "forecast"
$ws.Range("D28").Value = 46038
$ws.Range("E28").Value = "ON"
$ws.Range("A29").Value = 46038.5625
$ws.Range("B29").Value = 36.0601
$ws.Range("C29").Value = "forecast"
$ws.Range("D29").Value = 46038
$ws.Range("E29").Value = "ON"
$ws.Range("A30").Value = 46038.58333333334
$ws.Range("B30").Value = 36.0601
$ws.Range("C30").Value = "forecast"
$ws.Range("D30").Value = 46038
$ws.Range("E30").Value = "ON"
$ws.Range("A31").Value = 46038.60416666666
$ws.Range("B31").Value = 36.0601
$ws.Range("C31").Value = "forecast"
$ws.Range("D31").Value = 46038
$ws.Range("E31").Value = "ON"
$ws.Range("A32").Value = 46038.625
$ws.Range("B32").Value = 27.2786
$ws.Range("C32").Value = "forecast"
$ws.Range("D32").Value = 46038
$ws.Range("E32").Value = "ON"
$ws.Range("A33").Value = 46038.64583333334
$ws.Range("B33").Value = 18.70077
$ws.Range("C33").Value = "forecast"
$ws.Range("D33").Value = 46038
$ws.Range("E33").Value = "ON"
$ws.Range("A34").Value = 46038.66666666666
$ws.Range("B34").Value = 42.25471
$ws.Range("C34").Value = "forecast"
$ws.Range("D34").Value = 46038
$ws.Range("E34").Value = "OFF"
$ws.Range("A35").Value = 46038.6875
$ws.Range("B35").Value = 10.34966
$ws.Range("C35").Value = "forecast"
$ws.Range("D35").Value = 46038
$ws.Range("E35").Value = "OFF"
$ws.Range("A36").Value = 46038.70833333334
$ws.Range("B36").Value = 2.10968
$ws.Range("C36").Value = "forecast"
$ws.Range("D36").Value = 46038
$ws.Range("E36").Value = "OFF"
$ws.Range("A37").Value = 46038.72916666666
$ws.Range("B37").Value = -6
$ws.Range("C37").Value = "forecast"
$ws.Range("D37").Value = 46038
$ws.Range("E37").Value = "OFF"
$ws.Range("A38").Value = 46038.75
$ws.Range("B38").Value = -3.17523
$ws.Range("C38").Value = "forecast"
$ws.Range("D38").Value = 46038
$ws.Range("E38").Value = "OFF"
$ws.Range("A39").Value = 46038.77083333334
$ws.Range("B39").Value = -3.03165
$ws.Range("C39").Value = "forecast"
$ws.Range("D39").Value = 46038
$ws.Range("E39").Value = "OFF"
$ws.Range("A40").Value = 46038.79166666666
$ws.Range("B40").Value = 0.0113
$ws.Range("C40").Value = "forecast"
$ws.Range("D40").Value = 46038
$ws.Range("E40").Value = "OFF"
$ws.Range("A41").Value = 46038.8125
$ws.Range("B41").Value = 29.85322
$ws.Range("C41").Value = "forecast"
$ws.Range("D41").Value = 46038
$ws.Range("E41").Value = "OFF"
$ws.Range("A42").Value = 46038.83333333334
$ws.Range("B42").Value = 29.85322
$ws.Range("C42").Value = "forecast"
$ws.Range("D42").Value = 46038
$ws.Range("E42").Value = "ON"
$ws.Range("A43").Value = 46038.85416666666
$ws.Range("B43").Value = 29.85322
$ws.Range("C43").Value = "forecast"
$ws.Range("D43").Value = 46038
$ws.Range("E43").Value = "ON"
$ws.Range("A44").Value = 46038.875
$ws.Range("B44").Value = 8.67172
$ws.Range("C44").Value = "forecast"
$ws.Range("D44").Value = 46038
$ws.Range("E44").Value = "ON"
$ws.Range("A45").Value = 46038.89583333334
$ws.Range("B45").Value = 59.0817
$ws.Range("C45").Value = "forecast"
$ws.Range("D45").Value = 46038
$ws.Range("E45").Value = "ON"
$ws.Range("A46").Value = 46038.91666666666
$ws.Range("B46").Value = 57.09
$ws.Range("C46").Value = "forecast"
$ws.Range("D46").Value = 46038
$ws.Range("E46").Value = "ON"
$ws.Range("A47").Value = 46038.9375
$ws.Range("B47").Value = 57.44405
$ws.Range("C47").Value = "forecast"
$ws.Range("D47").Value = 46038
$ws.Range("E47").Value = "ON"
$ws.Range("A48").Value = 46038.95833333334
$ws.Range("B48").Value = 58.20349
$ws.Range("C48").Value = "forecast"
$ws.Range("D48").Value = 46038
$ws.Range("E48").Value = "ON"
$ws.Range("A49").Value = 46038.97916666666
$ws.Range("B49").Value = 56.98
$ws.Range("C49").Value = "forecast"
$ws.Range("D49").Value = 46038
$ws.Range("E49").Value = "ON"
